$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = "'29.952.91"
    'E2' = "'  -1.53%  "
    'D3' = "'1.893.27"
    'E3' = "'  -2.38%  "
    'D4' = "'1.000"
    'E4' = "'  -0.55%  "
    'D5' = "'0.7341"
    'E5' = "'  -1.69%  "
    'D6' = "'242.66"
    'E6' = "'  -1.18%  "
    'D7' = "'1.000"
    'E7' = "'  -0.52%  "
    'D8' = "'0.3096"
    'E8' = "'  -1.96%  "
    'E9' = "'  -4.53%  "
    'D10' = "'0.06897"
    'E10' = "'  -0.78%  "
    'D11' = "'0.7718"
    'E11' = "'  -0.91%  "
    'D12' = "'0.07953"
    'E12' = "'  -0.51%  "
    'D13' = "'1.890.97"
    'E13' = "'  -2.48%  "
    'D14' = "'5.221"
    'E14' = "'  -2.41%  "
    'D15' = "'91.53"
    'E15' = "'  -3.10%  "
    'B16' = "'WrappedBTC"
    'C16' = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
    'D16' = "'29.963.22"
    'E16' = "'  -1.54%  "
    'B17' = "'Avalanche"
    'C17' = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
    'D17' = "'14.13"
    'E17' = "'  -1.78%  "
    'D18' = "'5.784"
    'E18' = "'  +0.70%  "
    'D19' = "'240.06"
    'E19' = "'  -4.89%  "
    'D20' = "'0.000007749"
    'E20' = "'  -1.81%  "
    'D21' = "'1.0000"
    'E21' = "'  -0.44%  "
    'D22' = "'2.138.86"
    'E22' = "'  -2.37%  "
    'D23' = "'1.000"
    'E23' = "'  -0.60%  "
    'D24' = "'6.934"
    'E24' = "'  +4.03%  "
    'D25' = "'9.300"
    'E25' = "'  -1.89%  "
    'D26' = "'164.99"
    'E26' = "'  -0.36%  "
    'D27' = "'18.84"
    'E27' = "'  -0.62%  "
    'D28' = "'0.1270"
    'E28' = "'  -3.91%  "
    'E29' = "'  -10.05%  "
    'D30' = "'1.363"
    'E30' = "'  -0.13%  "
    'E31' = "'  +1.28%  "
    'D32' = "'4.307"
    'E32' = "'  -0.92%  "
    'D33' = "'4.060"
    'E33' = "'  -0.70%  "
    'D34' = "'0.05110"
    'E34' = "'  -0.69%  "
    'D35' = "'1.279"
    'E35' = "'  +0.52%  "
    'D36' = "'0.7360"
    'E36' = "'  -1.01%  "
    'D37' = "'2.717"
    'E37' = "'  -2.42%  "
    'D38' = "'0.01921"
    'E38' = "'  -1.02%  "
    'D39' = "'2.783"
    'E39' = "'  -0.78%  "
    'D40' = "'6.318"
    'E40' = "'  -1.56%  "
    'D41' = "'74.21"
    'E41' = "'  -4.91%  "
    'D42' = "'0.4458"
    'E42' = "'  +0.10%  "
    'D43' = "'1.934"
    'E43' = "'  -1.31%  "
    'D44' = "'1.001"
    'E44' = "'  -0.47%  "
    'D45' = "'0.8391"
    'E45' = "'  +0.87%  "
    'D46' = "'7.623"
    'E46' = "'  +2.39%  "
    'D47' = "'101.02"
    'E47' = "'  -0.24%  "
    'D48' = "'9.838"
    'E48' = "'  +0.65%  "
    'B49' = "'RocketPoolETH"
    'C49' = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
    'D49' = "'2.032.51"
    'E49' = "'  -2.55%  "
    'B50' = "'Elrond"
    'C50' = "'https://coinranking.com/coin/omwkOTglq+elrond-egld"
    'D50' = "'36.58"
    'E50' = "'  -1.56%  "
    'D51' = "'936.63"
    'E51' = "'  -4.03%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
